# Add a new "Player Info" worksheet as the first sheet in the workbook,
# and rework the existing "ODI Batting"/"ODI Bowling" sheets so the
# MATCH_CARD_LINK column becomes a bare MATCH_CODE value.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet before "ODI Batting" ---------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Style = "header"

# Data row
$playerInfo.Range("A2").Value = "4805"
$playerInfo.Range("B2").Value = "Neketh Gedara Roshan Prabath Jayasuriya"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# --- 2. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE -------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").Value = "4183"
$odiBatting.Range("D3").Value = "4186"

# --- 3. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE -------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"
$odiBowling.Range("B2").Value = "4183"
$odiBowling.Range("B3").Value = "4186"
